# Add an "index" loop variable to the list-row template.
#
# Original template rows (column A holds the template tag, the list loop
# spans rows 6-8, row 9 is the totals row):
#   A6: [row:list datalist as data]
#   A7: ${data.id}
#   A8: [/row:list]
#   A9: ${totalLabel}
#
# New template rows: the loop now also exposes an "index" variable and the
# per-row id cell prints that index instead of ${data.id}.
#   A6: [row:list datalist as data, index]
#   A7: ${index}
#   A9: ${totalLabel}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "[row:list datalist as data, index]"
$ws.Range("A7").Value = '${index}'
$ws.Range("A9").Value = '${totalLabel}'

# The author's last active selection moved down one row (A6:B6 -> A7:B7).
$ws.Range("A7:B7").Select()
